# Apply F-column (想去人数 / interested-count) updates across all 4 sheets
# per the commit 'Update gh-pages to output generated at 456a3b4'
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 86
$ws.Range("F4").Value = 9765
$ws.Range("F5").Value = 666
$ws.Range("F7").Value = 342
$ws.Range("F8").Value = 377
$ws.Range("F9").Value = 430
$ws.Range("F10").Value = 155
$ws.Range("F11").Value = 212
$ws.Range("F12").Value = 474
$ws.Range("F13").Value = 12404
$ws.Range("F14").Value = 38
$ws.Range("F19").Value = 251
$ws.Range("F21").Value = 183
$ws.Range("F22").Value = 127
$ws.Range("F24").Value = 2738
$ws.Range("F26").Value = 85
$ws.Range("F27").Value = 16
$ws.Range("F28").Value = 62
$ws.Range("F30").Value = 1048
$ws.Range("F31").Value = 4223
$ws.Range("F32").Value = 3713
$ws.Range("F33").Value = 704
$ws.Range("F35").Value = 3067
$ws.Range("F36").Value = 47
$ws.Range("F38").Value = 200
$ws.Range("F39").Value = 780
$ws.Range("F40").Value = 36
$ws.Range("F41").Value = 125
$ws.Range("F42").Value = 456
$ws.Range("F43").Value = 586
$ws.Range("F44").Value = 73
$ws.Range("F45").Value = 146
$ws.Range("F46").Value = 253
$ws.Range("F49").Value = 153

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F14").Value = 43

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 57

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 86
$ws.Range("F5").Value = 9765
$ws.Range("F6").Value = 666
$ws.Range("F9").Value = 342
$ws.Range("F10").Value = 377
$ws.Range("F11").Value = 430
$ws.Range("F12").Value = 155
$ws.Range("F13").Value = 212
$ws.Range("F14").Value = 474
$ws.Range("F15").Value = 12404
$ws.Range("F18").Value = 57
$ws.Range("F19").Value = 251
$ws.Range("F21").Value = 183
$ws.Range("F22").Value = 127
$ws.Range("F24").Value = 2738
$ws.Range("F26").Value = 85
$ws.Range("F27").Value = 62
$ws.Range("F29").Value = 1048
$ws.Range("F30").Value = 4223
$ws.Range("F31").Value = 3713
$ws.Range("F32").Value = 704
$ws.Range("F34").Value = 3067
$ws.Range("F35").Value = 47
$ws.Range("F37").Value = 200
$ws.Range("F38").Value = 780
$ws.Range("F39").Value = 36
$ws.Range("F40").Value = 125
$ws.Range("F41").Value = 456
$ws.Range("F43").Value = 587
$ws.Range("F44").Value = 73
$ws.Range("F45").Value = 146
$ws.Range("F46").Value = 253
$ws.Range("F49").Value = 153
